$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 203. This shifts every
# existing row from 203..269 down to 205..271 (and bumps the sheet
# dimension from A1:R269 to A1:R271), matching the rest of the diff, which
# is otherwise just every row's content sliding down by two positions.
$ws.Rows("203:204").Insert()

# Populate the two newly inserted rows with the new data points.

# Row 203 (new)
$ws.Cells.Item(203, 1).Value = 6
$ws.Cells.Item(203, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(203, 3).Value = "Metropolitana"
$ws.Cells.Item(203, 4).Value = 44468
$ws.Cells.Item(203, 5).Value = 13
$ws.Cells.Item(203, 6).Value = 100112043
$ws.Cells.Item(203, 7).Value = "Pepino ensalada"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 980
$ws.Cells.Item(203, 11).Value = 15000
$ws.Cells.Item(203, 12).Value = 16000
$ws.Cells.Item(203, 13).Value = 15459
$ws.Cells.Item(203, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(203, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(203, 16).Value = 258
$ws.Cells.Item(203, 17).Value = 60
$ws.Cells.Item(203, 18).Value = "Hortaliza"

# Row 204 (new)
$ws.Cells.Item(204, 1).Value = 6
$ws.Cells.Item(204, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(204, 3).Value = "Metropolitana"
$ws.Cells.Item(204, 4).Value = 44468
$ws.Cells.Item(204, 5).Value = 13
$ws.Cells.Item(204, 6).Value = 100112043
$ws.Cells.Item(204, 7).Value = "Pepino ensalada"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Segunda"
$ws.Cells.Item(204, 10).Value = 600
$ws.Cells.Item(204, 11).Value = 12000
$ws.Cells.Item(204, 12).Value = 14000
$ws.Cells.Item(204, 13).Value = 13167
$ws.Cells.Item(204, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(204, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(204, 16).Value = 132
$ws.Cells.Item(204, 17).Value = 100
$ws.Cells.Item(204, 18).Value = "Hortaliza"
